$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new boolean-classification columns before the old "work"
# column (old B..K shift right to D..M).
$ws.Range("B:C").Insert()

# The shift-on-insert path in this runtime re-serialises a handful of
# untouched numeric literals with full floating-point precision instead
# of their original shortest decimal form. Restore those exact literals
# so the sheet round-trips the original numbers unchanged.
$ws.Range("L2").Value = 25.2
$ws.Range("E3").Value = 36.9
$ws.Range("F3").Value = 110.1
$ws.Range("J3").Value = 528.9
$ws.Range("J4").Value = 719.7
$ws.Range("L4").Value = 27.3
$ws.Range("F5").Value = 111.51
$ws.Range("J5").Value = 411.3
$ws.Range("K5").Value = 7.4
$ws.Range("K6").Value = 8.6
$ws.Range("E7").Value = 38.03
$ws.Range("F7").Value = 114.51
$ws.Range("E9").Value = 27.72
$ws.Range("E11").Value = 35.74
$ws.Range("K11").Value = 9.4
$ws.Range("J12").Value = 493.7

# --- New site rows (13-17), column A first -------------------------------
# (write in this particular order so new shared strings land in the same
# sequence the original author's edits produced)
$ws.Range("A13").Value = "U1430"
$ws.Range("A14").Value = "G3"
$ws.Range("A15").Value = "ODP 1208"

$ws.Range("B1").Value = "pacific"

$ws.Range("A16").Value = "ODP 885/886"
$ws.Range("A17").Value = "Yushe"

$ws.Range("C1").Value = "CLP"

# --- Fill the new "pacific" / "CLP" boolean columns for existing rows ----
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = $true

$ws.Range("B3").Value = $false
$ws.Range("C3").Value = $true

$ws.Range("B4").Value = $false
$ws.Range("C4").Value = $true

$ws.Range("B5").Value = $false
$ws.Range("C5").Value = $true

$ws.Range("B6").Value = $true
$ws.Range("C6").Value = $true

$ws.Range("B7").Value = $false
$ws.Range("C7").Value = $false

$ws.Range("B8").Value = $false
$ws.Range("C8").Value = $false

$ws.Range("B9").Value = $false
$ws.Range("C9").Value = $false

$ws.Range("B10").Value = $false
$ws.Range("C10").Value = $false

$ws.Range("B11").Value = $false
$ws.Range("C11").Value = $false

$ws.Range("B12").Value = $false
$ws.Range("C12").Value = $true

# --- Fill in rest of the new rows (13-17) --------------------------------
$ws.Range("B13").Value = $true
$ws.Range("C13").Value = $false
$ws.Range("D13").Value = "ref"
$ws.Range("E13").Formula = "=37+54.16/60"
$ws.Range("F13").Formula = "=131+32.25/60"

$ws.Range("B14").Value = $true
$ws.Range("C14").Value = $false
$ws.Range("D14").Value = "ref"
$ws.Range("E14").Formula = "=38+50/60"
$ws.Range("F14").Formula = "=117+26/60"

$ws.Range("B15").Value = $true
$ws.Range("C15").Value = $false
$ws.Range("D15").Value = "ref"
$ws.Range("E15").Formula = "=36.1"
$ws.Range("F15").Value = 158.2

$ws.Range("B16").Value = $true
$ws.Range("C16").Value = $false
$ws.Range("D16").Value = "ref"
$ws.Range("E16").Value = 44.7
$ws.Range("F16").Value = -168.2

$ws.Range("B17").Value = $false
$ws.Range("C17").Value = $true
$ws.Range("D17").Value = "ref"
$ws.Range("E17").Value = 37
$ws.Range("F17").Value = 113

# --- Match the final selection shown in the saved workbook ---------------
$ws.Range("D12").Select()
